# Update "想去人数" (column F) figures across the four sheets to match the
# refreshed data snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 917
$ws.Range("F3").Value = 1482
$ws.Range("F4").Value = 1153
$ws.Range("F5").Value = 539
$ws.Range("F7").Value = 9
$ws.Range("F9").Value = 300
$ws.Range("F13").Value = 176
$ws.Range("F14").Value = 176
$ws.Range("F15").Value = 3959
$ws.Range("F16").Value = 24
$ws.Range("F17").Value = 18
$ws.Range("F18").Value = 448
$ws.Range("F22").Value = 411
$ws.Range("F23").Value = 118
$ws.Range("F26").Value = 70
$ws.Range("F27").Value = 273
$ws.Range("F28").Value = 977
$ws.Range("F30").Value = 1642
$ws.Range("F31").Value = 368

# --- 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 247
$ws.Range("F7").Value = 246
$ws.Range("F10").Value = 44

# --- 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 123

# --- 全部类型 (All types, combined view) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 917
$ws.Range("F4").Value = 1482
$ws.Range("F5").Value = 1153
$ws.Range("F8").Value = 123
$ws.Range("F9").Value = 539
$ws.Range("F11").Value = 9
$ws.Range("F14").Value = 300
$ws.Range("F18").Value = 176
$ws.Range("F19").Value = 176
$ws.Range("F20").Value = 3959
$ws.Range("F21").Value = 24
$ws.Range("F22").Value = 18
$ws.Range("F23").Value = 247
$ws.Range("F24").Value = 448
$ws.Range("F28").Value = 411
$ws.Range("F30").Value = 118
$ws.Range("F32").Value = 246
$ws.Range("F36").Value = 44
$ws.Range("F39").Value = 70
$ws.Range("F40").Value = 273
$ws.Range("F41").Value = 977
$ws.Range("F43").Value = 1642
$ws.Range("F44").Value = 368

$wb.Save()
